# Added Q7, which was missing from the BOM.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 19 (the R11... resistor row), shifting
# everything below it down by one. This is where the missing Q7 part
# (P-Channel MOSFET, same family as U8) belongs, alphabetically/numerically
# between Q6 (row 18) and the R-series resistors.
$ws.Rows.Item(19).Insert()

# Fill in the BOM data for the newly-added Q7 row.
$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "Q7"
$ws.Cells.Item(19, 3).Value = "Vishay"
$ws.Cells.Item(19, 4).Value = "SI7465DP-T1-E3"
$ws.Cells.Item(19, 5).Value = "P-Channel 60 V 3.2A (Ta) 1.5W (Ta) Surface Mount PowerPAK® SO-8"
$ws.Cells.Item(19, 6).Value = "Digi-Key"
$ws.Cells.Item(19, 7).Value = "SI7465DP-T1-E3CT-ND"

# Match the row's styling to the other highlighted rows (Q6 directly above,
# U8 which uses the same part) -- col C keeps the plain/default style while
# D/E pick up the wrapped variant of the highlight fill used elsewhere in
# this block.
$ws.Cells.Item(19, 3).Style = $ws.Cells.Item(18, 2).Style
$ws.Cells.Item(19, 4).Style = $ws.Cells.Item(18, 3).Style
$ws.Cells.Item(19, 5).Style = $ws.Cells.Item(18, 3).Style

# Move the active selection to reflect where the edit was made.
$ws.Range("D16").Select()
